$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for all touched cells (matches original inlineStr/text cells),
# then assign the new values per the diff.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.866.83"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.626.41"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.49"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.53"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.625.60"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.26%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.64%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.55"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.103.38"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.821.18"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.627.23"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.20"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.25%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.99"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.760.43"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.88%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "575.79"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.11%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.06%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.128"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.42%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "159.82"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.12"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.02%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.35"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.68%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "dogwifhat"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.03%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "BabyDogeCoin"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₆0334"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +17.00%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.05%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.10"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "154.90"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.88"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.21%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.33%  "
